# Applies the "Semana 8" (week 8) time-registry updates:
#  - Registro sheet: new duration entries for week 8 (rows 91-99 block)
#  - Total sheet: week-8 summary row (row 9) picks up the new Registro totals,
#    a running total formula is added in column H for rows 8-9, and the
#    grand-total row (row 10) recalculates accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Registro sheet - "Semana 8" block (rows 91-99)
# ---------------------------------------------------------------------
$reg = $wb.Worksheets.Item("Registro")

# New time entries that previously were blank.
$reg.Range("B92").Value2 = 0.020833333333333332
$reg.Range("B92").NumberFormat = "h:mm:ss"
$reg.Range("B92").HorizontalAlignment = 1
$reg.Range("B92").Font.Name = "Arial"

$reg.Range("F92").Value2 = 0.03125
$reg.Range("F92").NumberFormat = "h:mm:ss"
$reg.Range("F92").HorizontalAlignment = 1
$reg.Range("F92").Font.Name = "Arial"

$reg.Range("G92").Value2 = 0.041666666666666664
$reg.Range("G92").NumberFormat = "h:mm:ss"
$reg.Range("G92").HorizontalAlignment = 1
$reg.Range("G92").Font.Name = "Arial"

# Existing entry updated with a larger duration.
$reg.Range("C95").Value2 = 0.24097222222222223

# New entry in a cell that previously had no value.
$reg.Range("C96").Value2 = 0.10609953703703703
$reg.Range("C96").NumberFormat = "h:mm:ss"
$reg.Range("C96").HorizontalAlignment = -4152
$reg.Range("C96").Font.Name = "Arial"

# ---------------------------------------------------------------------
# 2) Total sheet - week 8 row (row 9) mirrors the new Registro totals
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("Total")

$total.Range("B9").Value2 = 0.06805555555555555
$total.Range("B9").NumberFormat = "h:mm:ss"

$total.Range("C9").Value2 = 0.7815856481481481
$total.Range("F9").Value2 = 0.2222222222222222
$total.Range("G9").Value2 = 0.11805555555555555

# Running per-week total in column H (previously blank for rows 8 and 9),
# entered as one shared formula across H8:H9.
$total.Range("H8:H9").Formula = "=SUM(B8:G8)"

# Row 10 (grand totals) and H8/H9 recompute automatically through the
# existing SUM formulas once the source cells above are updated.
